# Documentation update: add a "Simulation" column (after Trial) and a
# "Comments" column (at the end), and correct a couple of "Done" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Obstacles" column (B) to hold
# the new "Simulation" data. This shifts B:L -> C:M, preserving all
# existing values/styles in their new positions.
$ws.Columns("B").Insert()

# --- Header row ---
$ws.Range("B1").Value = "Simulation"
$ws.Range("N1").Value = "Comments"

# --- New "Simulation" column values (all trials were run in simulation) ---
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1

# --- Corrected "Done" values (column M after the insert) ---
$ws.Range("M3").Value = "Yes"
$ws.Range("M6").Value = "Yes"

# --- New "Comments" column values ---
$ws.Range("N2").Value = "Bad"
$ws.Range("N3").Value = "Bad"
$ws.Range("N5").Value = "Very Bad"
$ws.Range("N6").Value = "Ok"
$ws.Range("N7").Value = "Meh"
$ws.Range("N8").Value = "Bad"

# --- Update the selection to span the full used range ---
$ws.Range("A1:N8").Select()
